$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.619.91'
$ws.Range('E2').Value = '  +2.73%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.479.56'
$ws.Range('E3').Value = '  +2.53%  '

# Row 4
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '574.90'
$ws.Range('E5').Value = '  +2.02%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.19'
$ws.Range('E6').Value = '  +5.15%  '

# Row 7
$ws.Range('E7').Value = '  -0.13%  '

# Row 8
$ws.Range('E8').Value = '  +2.03%  '

# Row 9
$ws.Range('E9').Value = '  +4.60%  '

# Row 10
$ws.Range('E10').Value = '  +0.29%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.365'
$ws.Range('E11').Value = '  +4.41%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.34'
$ws.Range('E12').Value = '  +2.57%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.36'
$ws.Range('E13').Value = '  +6.06%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000185'
$ws.Range('E14').Value = '  +7.50%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.900.72'
$ws.Range('E15').Value = '  +1.61%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.415.22'
$ws.Range('E16').Value = '  +2.52%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.495.06'
$ws.Range('E17').Value = '  +3.06%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.61'
$ws.Range('E18').Value = '  +2.83%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.25'
$ws.Range('E19').Value = '  +6.34%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.26'
$ws.Range('E20').Value = '  +3.51%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '329.42'
$ws.Range('E21').Value = '  +2.01%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.11%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.90'
$ws.Range('E23').Value = '  +10.03%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '67.68'
$ws.Range('E24').Value = '  +1.66%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '643.41'
$ws.Range('E25').Value = '  +15.62%  '

# Row 26
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.84'
$ws.Range('E26').Value = '  +0.68%  '

# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  +13.16%  '

# Row 28
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.607.30'
$ws.Range('E28').Value = '  +2.77%  '

# Row 29
$ws.Range('B29').Value = 'Fetch.AI'
$ws.Range('C29').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.53'
$ws.Range('E29').Value = '  +10.49%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.55'
$ws.Range('E30').Value = '  +4.59%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.985'
$ws.Range('E31').Value = '  -1.51%  '

# Row 32
$ws.Range('E32').Value = '  -1.58%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.92'
$ws.Range('E33').Value = '  +2.91%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.25'
$ws.Range('E34').Value = '  +10.68%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.57'
$ws.Range('E35').Value = '  +4.82%  '

# Row 36
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.16%  '

# Row 37
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.388'
$ws.Range('E37').Value = '  +2.51%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.53'
$ws.Range('E38').Value = '  +1.96%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.00'
$ws.Range('E39').Value = '  +2.72%  '

# Row 40
$ws.Range('E40').Value = '  +2.93%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '147.47'
$ws.Range('E41').Value = '  -4.13%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.63'
$ws.Range('E42').Value = '  +17.37%  '

# Row 43
$ws.Range('E43').Value = '  +0.60%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '151.91'
$ws.Range('E44').Value = '  +3.28%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.81'
$ws.Range('E45').Value = '  +4.85%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0555'
$ws.Range('E46').Value = '  +5.67%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '21.17'
$ws.Range('E47').Value = '  +7.11%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.612'
$ws.Range('E48').Value = '  +3.41%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0241'
$ws.Range('E49').Value = '  +6.17%  '

# Row 50
$ws.Range('E50').Value = '  +1.12%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.741'
$ws.Range('E51').Value = '  +5.12%  '
